# Trans-Tasman charts updated to 29 November
# Appends the latest rows of NZ vaccination-by-date data (through 28 Nov 2021,
# serial 44528) to Sheet1, mirroring the existing table's layout/format, and
# updates the sheet's active selection to reflect the new last row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New daily rows to append after the existing last row (278).
# Columns: A = date (serial), B = first doses, C = second doses
$newRows = @(
  @(44522, 4777, 10447),
  @(44523, 6496, 12384),
  @(44524, 6600, 12023),
  @(44525, 6319, 11773),
  @(44526, 5679, 12222),
  @(44527, 6378, 14495),
  @(44528, 3679, 8040)
)

$lastRow = 278
$startRow = $lastRow + 1

for ($i = 0; $i -lt $newRows.Length; $i++) {
  $r = $startRow + $i

  # Copy the formatting of the last existing row down into the new row
  # (keeps the date number format / alignment style used by column A).
  $ws.Range("A" + $lastRow + ":C" + $lastRow).Copy($ws.Range("A" + $r + ":C" + $r))

  $row = $newRows[$i]
  $ws.Range("A" + $r).Value = $row[0]
  $ws.Range("B" + $r).Value = $row[1]
  $ws.Range("C" + $r).Value = $row[2]
}

$newLastRow = $startRow + $newRows.Length - 1

# Scroll/select to match where the author left off: the view is scrolled to
# row 259, and the previously-selected cell (old last row, column C) is kept
# selected, which now sits at C278 after the 7 new rows were appended.
$win = $excel.ActiveWindow
$win.ScrollRow = 259
$win.ScrollColumn = 1
$ws.Range("C" + $lastRow).Select()
